$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 988 (pushes existing rows 988-1034 down to 991-1037)
$ws.Rows.Item(988).EntireRow.Insert()
$ws.Rows.Item(988).EntireRow.Insert()
$ws.Rows.Item(988).EntireRow.Insert()

# New weekly price entries (fecha 45075) for Choclo Lluteño - Primera/Segunda/Tercera
$ws.Cells.Item(988, 1).Value = 1
$ws.Cells.Item(988, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(988, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(988, 4).Value = 45075
$ws.Cells.Item(988, 5).Value = 15
$ws.Cells.Item(988, 6).Value = 100112024
$ws.Cells.Item(988, 7).Value = "Choclo"
$ws.Cells.Item(988, 8).Value = "Lluteño"
$ws.Cells.Item(988, 9).Value = "Primera"
$ws.Cells.Item(988, 10).Value = 40
$ws.Cells.Item(988, 11).Value = 34000
$ws.Cells.Item(988, 12).Value = 35000
$ws.Cells.Item(988, 13).Value = 34500
$ws.Cells.Item(988, 14).Value = "$/saco 50 unidades"
$ws.Cells.Item(988, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(988, 16).Value = 690
$ws.Cells.Item(988, 17).Value = 50
$ws.Cells.Item(988, 18).Value = "Hortaliza"

$ws.Cells.Item(989, 1).Value = 1
$ws.Cells.Item(989, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(989, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(989, 4).Value = 45075
$ws.Cells.Item(989, 5).Value = 15
$ws.Cells.Item(989, 6).Value = 100112024
$ws.Cells.Item(989, 7).Value = "Choclo"
$ws.Cells.Item(989, 8).Value = "Lluteño"
$ws.Cells.Item(989, 9).Value = "Segunda"
$ws.Cells.Item(989, 10).Value = 40
$ws.Cells.Item(989, 11).Value = 29000
$ws.Cells.Item(989, 12).Value = 30000
$ws.Cells.Item(989, 13).Value = 29500
$ws.Cells.Item(989, 14).Value = "$/saco 75 unidades"
$ws.Cells.Item(989, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(989, 16).Value = 393
$ws.Cells.Item(989, 17).Value = 75
$ws.Cells.Item(989, 18).Value = "Hortaliza"

$ws.Cells.Item(990, 1).Value = 1
$ws.Cells.Item(990, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(990, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(990, 4).Value = 45075
$ws.Cells.Item(990, 5).Value = 15
$ws.Cells.Item(990, 6).Value = 100112024
$ws.Cells.Item(990, 7).Value = "Choclo"
$ws.Cells.Item(990, 8).Value = "Lluteño"
$ws.Cells.Item(990, 9).Value = "Tercera"
$ws.Cells.Item(990, 10).Value = 50
$ws.Cells.Item(990, 11).Value = 24000
$ws.Cells.Item(990, 12).Value = 25000
$ws.Cells.Item(990, 13).Value = 24500
$ws.Cells.Item(990, 14).Value = "$/saco 100 unidades"
$ws.Cells.Item(990, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(990, 16).Value = 245
$ws.Cells.Item(990, 17).Value = 100
$ws.Cells.Item(990, 18).Value = "Hortaliza"
